# Add two new columns: I ("I0") and J ("IF")
# I0 is always 1, IF mirrors the existing IP (column H) value for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy the formatting of the existing "IP" header (H1)
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

$excel.CutCopyMode = 0

# Data rows 2 through 31
for ($r = 2; $r -le 31; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2  # column H = "IP"

    $ws.Cells.Item($r, 9).Value2 = 1          # column I = "I0"
    $ws.Cells.Item($r, 10).Value2 = $ipValue  # column J = "IF" (same as IP)
}
